$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.477.51"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "3.435.68"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'412.93"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "'128.19"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.731"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'42.89"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "'0.0000220"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "3.987.52"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'20.57"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "3.444.83"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'12.73"
$ws.Range("E18").Value = "  +5.01%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "62.360.30"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'479.70"
$ws.Range("E21").Value = "  +9.01%  "
$ws.Range("D22").Value = "'91.96"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "'3.28"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").Value = "'13.07"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'3.31"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "'9.67"
$ws.Range("E26").Value = "  +8.75%  "
$ws.Range("D27").Value = "'33.58"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'4.79"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.75"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'11.88"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").Value = "'41.08"
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'58.08"
$ws.Range("E36").Value = "  +7.78%  "
$ws.Range("D37").Value = "'0.0490"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.09"
$ws.Range("E38").Value = "  +6.09%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'148.35"
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.324"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "'2.08"
$ws.Range("E44").Value = "  +4.73%  "
$ws.Range("D45").Value = "'2.61"
$ws.Range("E45").Value = "  +8.40%  "
$ws.Range("D46").Value = "'4.26"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("D47").Value = "0.0₃0551"
$ws.Range("E47").Value = "  +26.66%  "
$ws.Range("B48").Value = "Celestia"
$ws.Range("C48").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D48").Value = "'16.39"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.30"
$ws.Range("E49").Value = "  +16.16%  "
$ws.Range("D50").Value = "'22.29"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'114.61"
$ws.Range("E51").Value = "  +8.82%  "
